$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.757810711860657
$ws.Range("B1").Value = 2.686997413635254
$ws.Range("C1").Value = 3.338043928146362
$ws.Range("D1").Value = 1.27691650390625
$ws.Range("E1").Value = 0.8488638997077942
